$wb = $excel.ActiveWorkbook

$wsTemplate = $wb.Worksheets.Item("Template")
$wsPrint1 = $wb.Worksheets.Item("print 1")

# --- "print 1" sheet: shipping/customer details + items ---
$wsPrint1.Range("A2").Value = "Asia"
$wsPrint1.Range("A3").Value = "Jl. Cinere Raya No.202"
$wsPrint1.Range("A4").Value = "Cinere, Kota Depok"
$wsPrint1.Range("A5").Value = "Jawa Barat, 16514"

# --- "Template" sheet: reprint marker + customer label ---
$wsTemplate.Range("G5").Value = "[reprint]"
$wsTemplate.Range("A1").Value = "Customer :[count]"

# --- "print 1" sheet: item rows ---
$wsPrint1.Range("A7").Value = "TFL-138-22"
$wsPrint1.Range("C7").Value = "Box"
$wsPrint1.Range("A8").Value = "TA-1006-KB"
$wsPrint1.Range("C8").Value = "Box"

# --- "print 1" sheet: print timestamp + page/print count ---
$wsPrint1.Range("G2").Value = "02-10-2024 06:15"
$wsPrint1.Range("G4").Value = 117

# --- "print 1" sheet: quantities ---
$wsPrint1.Range("B7").Value = 4
$wsPrint1.Range("B8").Value = 7

# Reset the lingering cell selection left over from editing (both sheets
# were last parked on A5); move back to the top-left "home" cell and make
# sure "Template" ends up the active tab again, matching the saved view.
$wsPrint1.Range("A1").Select()
$wsTemplate.Activate()
$wsTemplate.Range("A1").Select()
